$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.042494666666667
$ws.Range("H2").Value = 3.127484
$ws.Range("I2").Value = 0.0007670466909205676
$ws.Range("J2").Value = 0.0007670466909205677
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2615913333333333
$ws.Range("N2").Value = 0.784774
$ws.Range("O2").Value = 0.08239613548481725
$ws.Range("P2").Value = 0.08239613548481727
$ws.Range("Q2").Value = 0.2727075698462222
$ws.Range("R2").Value = 2.454368128616
$ws.Range("S2").Value = [double]"6.320168306827183E-05"
$ws.Range("T2").Value = [double]"6.320168306827185E-05"

$ws.Range("G3").Value = 1.042494666666667
$ws.Range("H3").Value = 3.127484
$ws.Range("I3").Value = 0.0007670466909205676
$ws.Range("J3").Value = 0.0007670466909205677
$ws.Range("N3").Value = 5.233242000000001
$ws.Range("O3").Value = 0.5494561706387266
$ws.Range("P3").Value = 0.5494561706387268
$ws.Range("Q3").Value = 1.818542291458667
$ws.Range("R3").Value = 16.366880623128
$ws.Range("S3").Value = 0.000421458537494322
$ws.Range("T3").Value = 0.0004214585374943222

$ws.Range("G4").Value = 1.042494666666667
$ws.Range("H4").Value = 3.127484
$ws.Range("I4").Value = 0.0007670466909205676
$ws.Range("J4").Value = 0.0007670466909205677
$ws.Range("M4").Value = 1.168795666666667
$ws.Range("N4").Value = 3.506387
$ws.Range("O4").Value = 0.3681476938764561
$ws.Range("P4").Value = 0.3681476938764561
$ws.Range("Q4").Value = 1.218463248923111
$ws.Range("R4").Value = 10.966169240308
$ws.Range("S4").Value = 0.0002823864703579737
$ws.Range("T4").Value = 0.0002823864703579738

$ws.Range("I5").Value = 0.9658609009611662
$ws.Range("J5").Value = 0.9658609009611662
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2615913333333333
$ws.Range("N5").Value = 0.784774
$ws.Range("O5").Value = 0.08239613548481725
$ws.Range("P5").Value = 0.08239613548481727
$ws.Range("Q5").Value = 343.3918459311609
$ws.Range("R5").Value = 3090.526613380448
$ws.Range("S5").Value = 0.07958320565508391
$ws.Range("T5").Value = 0.07958320565508392

$ws.Range("I6").Value = 0.9658609009611662
$ws.Range("J6").Value = 0.9658609009611662
$ws.Range("N6").Value = 5.233242000000001
$ws.Range("O6").Value = 0.5494561706387266
$ws.Range("P6").Value = 0.5494561706387268
$ws.Range("Q6").Value = 2289.898277191243
$ws.Range("S6").Value = 0.5306982320117928
$ws.Range("T6").Value = 0.5306982320117929

$ws.Range("I7").Value = 0.9658609009611662
$ws.Range("J7").Value = 0.9658609009611662
$ws.Range("M7").Value = 1.168795666666667
$ws.Range("N7").Value = 3.506387
$ws.Range("O7").Value = 0.3681476938764561
$ws.Range("P7").Value = 0.3681476938764561
$ws.Range("Q7").Value = 1534.282104757581
$ws.Range("R7").Value = 13808.53894281822
$ws.Range("S7").Value = 0.3555794632942895
$ws.Range("T7").Value = 0.3555794632942895

$ws.Range("G8").Value = 45.356022
$ws.Range("H8").Value = 136.068066
$ws.Range("I8").Value = 0.03337205234791334
$ws.Range("J8").Value = 0.03337205234791334
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2615913333333333
$ws.Range("N8").Value = 0.784774
$ws.Range("O8").Value = 0.08239613548481725
$ws.Range("P8").Value = 0.08239613548481727
$ws.Range("Q8").Value = 11.864742269676
$ws.Range("R8").Value = 106.782680427084
$ws.Range("S8").Value = 0.002749728146665081
$ws.Range("T8").Value = 0.002749728146665081

$ws.Range("G9").Value = 45.356022
$ws.Range("H9").Value = 136.068066
$ws.Range("I9").Value = 0.03337205234791334
$ws.Range("J9").Value = 0.03337205234791334
$ws.Range("N9").Value = 5.233242000000001
$ws.Range("O9").Value = 0.5494561706387266
$ws.Range("P9").Value = 0.5494561706387268
$ws.Range("Q9").Value = 79.119679761108
$ws.Range("R9").Value = 712.0771178499721
$ws.Range("S9").Value = 0.01833648008943959
$ws.Range("T9").Value = 0.01833648008943959

$ws.Range("G10").Value = 45.356022
$ws.Range("H10").Value = 136.068066
$ws.Range("I10").Value = 0.03337205234791334
$ws.Range("J10").Value = 0.03337205234791334
$ws.Range("M10").Value = 1.168795666666667
$ws.Range("N10").Value = 3.506387
$ws.Range("O10").Value = 0.3681476938764561
$ws.Range("P10").Value = 0.3681476938764561
$ws.Range("Q10").Value = 53.011921970838
$ws.Range("R10").Value = 477.107297737542
$ws.Range("S10").Value = 0.01228584411180867
$ws.Range("T10").Value = 0.01228584411180867
